$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free direct cell updates mirroring the refreshed coinranking.com scrape.
# Numeric-looking price strings (e.g. '1.00', '61.00') are forced to Text so Excel
# doesn't collapse the trailing zeros / decimal point into a real number; the cell's
# original (unstyled) format is restored immediately after so no stray style sticks.

$ws.Range("D2").Value = '26.214.58'
$ws.Range("E2").Value = '  -0.41%  '

$ws.Range("D3").Value = '1.586.54'
$ws.Range("E3").Value = '  -0.27%  '

$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.98%  '

$ws.Range("E6").Value = '  -0.43%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.15%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.245'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.23%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0604'
$ws.Range("D9").Style = "Normal"

$ws.Range("E10").Value = '  -1.70%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0846'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.18%  '

$ws.Range("D12").Value = '1.809.29'

$ws.Range("D13").Value = '1.598.61'
$ws.Range("E13").Value = '  +0.62%  '

$ws.Range("E14").Value = '  -1.55%  '

$ws.Range("E15").Value = '  -0.41%  '

$ws.Range("E16").Value = '  -0.87%  '

$ws.Range("D17").Value = '26.210.58'
$ws.Range("E17").Value = '  -0.43%  '

$ws.Range("D18").Value = '0.0₃0725'
$ws.Range("E18").Value = '  -0.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '214.04'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.53%  '

$ws.Range("E21").Value = '  -0.09%  '

$ws.Range("E23").Value = '  +0.46%  '

$ws.Range("E24").Value = '  -1.71%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.68'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.15%  '

$ws.Range("E26").Value = '  -0.16%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.84%  '

$ws.Range("E28").Value = '  -0.91%  '

$ws.Range("E29").Value = '  -1.22%  '

$ws.Range("E30").Value = '  -1.85%  '

$ws.Range("E31").Value = '  +0.70%  '

$ws.Range("E32").Value = '  -0.83%  '

$ws.Range("D33").Value = '1.416.59'
$ws.Range("E33").Value = '  +8.41%  '

$ws.Range("E34").Value = '  -1.66%  '

$ws.Range("E35").Value = '  -0.59%  '

$ws.Range("E36").Value = '  -1.49%  '

$ws.Range("E37").Value = '  -4.74%  '

$ws.Range("E38").Value = '  -1.28%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.91'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.22%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.822'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.24%  '

$ws.Range("E41").Value = '  -0.21%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.953'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -12.79%  '

$ws.Range("E43").Value = '  +0.24%  '

$ws.Range("E44").Value = '  -0.08%  '

$ws.Range("D45").Value = '1.720.88'
$ws.Range("E45").Value = '  -0.27%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.48%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '85.54'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.77%  '

$ws.Range("E48").Value = '  -0.85%  '

$ws.Range("E49").Value = '  -0.70%  '

$ws.Range("E50").Value = '  -1.35%  '

$ws.Range("E51").Value = '  -0.19%  '
